$wb = $excel.ActiveWorkbook
$sheets = @($wb.Worksheets.Item("展览"), $wb.Worksheets.Item("全部类型"))

foreach ($ws in $sheets) {
    # --- Simple "want to go" count bumps in column F ---
    $ws.Cells.Item(5, 6).Value = 3074
    $ws.Cells.Item(7, 6).Value = 2401
    $ws.Cells.Item(10, 6).Value = 2
    $ws.Cells.Item(11, 6).Value = 1195
    $ws.Cells.Item(14, 6).Value = 6
    $ws.Cells.Item(15, 6).Value = 1075
    $ws.Cells.Item(16, 6).Value = 289
    $ws.Cells.Item(17, 6).Value = 309
    $ws.Cells.Item(20, 6).Value = 102
    $ws.Cells.Item(21, 6).Value = 60

    # --- Insert a new row after row 23 ---
    # (old row 23 "南昌·代号鸢盛花行only" shifts to row 24,
    #  old row 24 "九江·第三届ACD动漫游戏嘉年华" shifts to row 25)
    $ws.Rows.Item(24).Insert()

    # Copy the row-number cell's formatting down so A24 matches A22/A23's style
    $ws.Range("A22").Copy()
    $ws.Range("A24").PasteSpecial(-4122)

    # Row 24 now holds what used to be row 23's content
    $ws.Cells.Item(24, 1).Value = 23
    $ws.Cells.Item(24, 2).Value = "'2024-04-20"
    $ws.Cells.Item(24, 3).Value = "南昌·代号鸢盛花行only"
    $ws.Cells.Item(24, 4).Value = "民德路411号 东方豪景花园酒店(民德路店)"
    $ws.Cells.Item(24, 5).Value = "2024.04.20 09:30-04.20 17:30"
    $ws.Cells.Item(24, 6).Value = 5
    $ws.Cells.Item(24, 7).Value = 78
    $ws.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82529"
    $ws.Cells.Item(24, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png"

    # The quoted-date entry above leaves a "number stored as text" style on
    # B24; clear it by repainting with B23's (unstyled) formatting.
    $ws.Range("B23").Copy()
    $ws.Range("B24").PasteSpecial(-4122)

    # Row 23 becomes the brand-new "New World" event
    $ws.Cells.Item(23, 3).Value = "南昌·New World国潮动漫博览会"
    $ws.Cells.Item(23, 4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Cells.Item(23, 5).Value = "2024.04.20 09:30-04.21 17:00"
    $ws.Cells.Item(23, 6).Value = 15
    $ws.Cells.Item(23, 7).Value = 60
    $ws.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82411"
    $ws.Cells.Item(23, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/xbYbLXc81709707724935.jpeg"

    # Row 25 (old row 24, shifted down): bump its index 23 -> 24 and F count 31 -> 32
    $ws.Cells.Item(25, 1).Value = 24
    $ws.Cells.Item(25, 6).Value = 32
}
